$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SA")

# Rows 1-9 (column headers, HKL/reflection headers, and the first 7 existing
# schemes: ND Single .. Ring Perpendicular to TD) are untouched by this
# change, so we leave them exactly as they are.

# Rows 10-16: the scheme list is reordered so that "Gaussian-Quadrature"
# (previously the last scheme, on row 16) now comes right after "Ring
# Perpendicular to TD", followed by three brand-new spiral schemes, and then
# the schemes that used to start at row 10 continue on (shifted down by 4).
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"

# Rows 17-19: brand new rows holding the remaining pre-existing schemes that
# got pushed further down (HexGrid-*). Give column A the same
# bold/centered/bordered style used by the rest of the index column (copied
# from A16 so styles.xml gains no new entries), and fill in the averaged
# intensity values (1) across C:P like every other data row.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"

for ($row = 17; $row -le 19; $row++) {
  for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item($row, $col).Value = 1
  }
}

$ws.Range("A1").Select() | Out-Null
